$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for column D
$ws.Range("D1").Value = "Corrected_Splitting"

# New data values for column D, rows 2-12
$values = @(103, 105, 103.4, 104, 104.5, 104.1, 105.2, 103.5, 103, 103.4, 104)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Update the selection to match the diff (D13)
$ws.Range("D13").Select()
